# task 3 and 5
# Re-highlight (white -> yellow) the text runs belonging to the "Реализовать
# воспроизведение видеофайла..." list item (task 3) and most of the
# "Создать элемент на странице..." list item (task 5). For task 5 the final
# sentence is split: only "закрываться. " becomes yellow, the remainder
# ("Необходимо использовать ...дебаунс.") stays white.

$d = $word.ActiveDocument

function Set-HighlightOnText($text) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $r.Find.Replacement.ClearFormatting()
    $r.Find.Text = $text
    $r.Find.Replacement.Text = $text
    $r.Find.Replacement.Highlight = $true
    # wdReplaceOne = 1 -- replace just the (unique) first/only match we target
    $r.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 1) | Out-Null
}

# --- Task 3 paragraph: "Реализовать воспроизведение видеофайла ..." ---
Set-HighlightOnText "Реализовать воспроизведение видеофайла (желательно формата "
Set-HighlightOnText "webm"
Set-HighlightOnText "). Никаких элементов управления на блоке быть не должно, только "
Set-HighlightOnText "предзагруженный"
Set-HighlightOnText " первый кадр видео и текущее время видео в формате "
Set-HighlightOnText "MM:SS"
Set-HighlightOnText ":mmm"
Set-HighlightOnText ". Видео должно запускаться и ставиться на паузу по клику мыши. По окончанию видео останавливается и возвращается на начальный кадр;"

# --- Task 5 paragraphs: "Создать элемент на странице ..." ---
Set-HighlightOnText "Создать элемент на странице. "
Set-HighlightOnText "При клике"
Set-HighlightOnText " на который будут открываться 3 любых"
Set-HighlightOnText "элемента на странице в течении 1 секунды. И при повторном клике 3 элемента будут"

# Only the "закрываться. " portion of the run turns yellow; this also
# splits the original run so "Необходимо использовать " remains its own
# (still white) run.
Set-HighlightOnText "закрываться. "
